# MongoDB索引.pptx edit script
# 1. Expand slide 3 ("MongoDB有哪些索引") content with 2 extra bullets (Hash索引, 全文索引)
# 2. Fill in slide 5 ("覆盖索引") content (previously empty)
# 3. Insert 4 new slides after slide 3: "_id索引", "单字段索引", "复合索引", "多键索引"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Step 1: slide 3 "MongoDB有哪些索引" - insert "Hash索引" and "全文索引" paragraphs
# right before the existing "过期索引" paragraph.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

$body3.Text = "_id索引`r单字段索引`r复合索引`r多键索引`rHash索引`r全文索引`r过期索引`r`r"

# fix up the "Hash" run of paragraph 5 so "Hash" and "索引" are separate runs
$para5 = $body3.Paragraphs(5)
$para5.Text = "索引"
$hashRun = $para5.InsertBefore("Hash")
$hashRun.LanguageID = "en-US"
$para5.LanguageID = "zh-CN"

$para1 = $body3.Paragraphs(1)
$para1.Text = "索引"
$idRun = $para1.InsertBefore("_id")
$idRun.LanguageID = "en-US"
$para1.LanguageID = "zh-CN"

# ---------------------------------------------------------------------------
# Step 2: slide 5 "覆盖索引" - fill in the previously-empty content placeholder
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange

$body5.Text = "当查询符合以下条件时，MongoDB可以直接从索引中获取返回数据，而不需要回表查询整个集合。`r1. 所有的查询字段是索引的一部分`r2. 所有的查询返回字段在同一个索引中`r因为索引存在于RAM中，从索引中获取数据比通过扫描文档读取数据要快得多。`r创建的索引中一般不包括 _id 字段，而_id在查询中会默认返回，我们可以在MongoDB的查询结果集中排除_id字段来实现索引覆盖。"

$p1 = $body5.Paragraphs(1)
$p1.Text = "可以直接从索引中获取返回数据，而不需要回表查询整个集合。"
$mongoRun1 = $p1.InsertBefore("MongoDB")
$mongoRun1.LanguageID = "en-US"
$beforeRun1 = $mongoRun1.InsertBefore("当查询符合以下条件时，")
$beforeRun1.LanguageID = "zh-CN"
$p1.LanguageID = "zh-CN"

$p2 = $body5.Paragraphs(2)
$p2.Text = " 所有的查询字段是索引的一部分"
$num1 = $p2.InsertBefore("1.")
$num1.LanguageID = "en-US"
$p2.LanguageID = "zh-CN"
$p2.ParagraphFormat.Bullet.Visible = 0
$p2.IndentLevel = 1

$p3 = $body5.Paragraphs(3)
$p3.Text = " 所有的查询返回字段在同一个索引中"
$num2 = $p3.InsertBefore("2.")
$num2.LanguageID = "en-US"
$p3.LanguageID = "zh-CN"
$p3.ParagraphFormat.Bullet.Visible = 0
$p3.IndentLevel = 1

$p4 = $body5.Paragraphs(4)
$p4.Text = "中，从索引中获取数据比通过扫描文档读取数据要快得多。"
$ramRun = $p4.InsertBefore("RAM")
$ramRun.LanguageID = "en-US"
$beforeRun4 = $ramRun.InsertBefore("因为索引存在于")
$beforeRun4.LanguageID = "zh-CN"
$p4.LanguageID = "zh-CN"

$p5 = $body5.Paragraphs(5)
$p5.Text = "字段来实现索引覆盖。"
$r1 = $p5.InsertBefore("创建的索引中一般不包括 ")
$r1.LanguageID = "zh-CN"
$r2 = $r1.InsertAfter("_")
$r2.LanguageID = "en-US"
$r3 = $r2.InsertAfter("id ")
$r3.LanguageID = "en-US"
$r4 = $r3.InsertAfter("字段，而")
$r4.LanguageID = "zh-CN"
$r5 = $r4.InsertAfter("_")
$r5.LanguageID = "en-US"
$r6 = $r5.InsertAfter("id")
$r6.LanguageID = "en-US"
$r7 = $r6.InsertAfter("在查询中会默认返回，我们可以在")
$r7.LanguageID = "zh-CN"
$r8 = $r7.InsertAfter("MongoDB")
$r8.LanguageID = "en-US"
$r9 = $r8.InsertAfter("的查询结果集中排除")
$r9.LanguageID = "zh-CN"
$r10 = $r9.InsertAfter("_id")
$r10.LanguageID = "en-US"
$p5.LanguageID = "zh-CN"

# ---------------------------------------------------------------------------
# Step 3: insert 4 new slides after slide 3 (index 4,5,6,7):
#   "_id索引", "单字段索引", "复合索引", "多键索引"
# Use layout index 2 (Title and Content) to match other content slides.
# ---------------------------------------------------------------------------

# --- New slide 4: "_id索引" ---
$nsId = $p.Slides.Add(4, 2)
$nsId.Shapes.Item(1).Name = "标题 1"
$nsId.Shapes.Item(2).Name = "内容占位符 2"

$titleId = $nsId.Shapes.Item(1).TextFrame.TextRange
$titleId.Text = "索引"
$titleRun = $titleId.InsertBefore("_id")
$titleRun.LanguageID = "en-US"
$titleId.LanguageID = "zh-CN"

$bodyId = $nsId.Shapes.Item(2).TextFrame.TextRange
$bodyId.Text = "属于单字段索引的一种，MongoDB默认创建`r默认情况下，_id 字段的类型为 ObjectID，是 MongoDB 的 BSON 类型之一。`rObjectID 长度为 12 字节，由以下4个部分组成：`r4 字节的Unix时间戳，单位为秒`r3 字节的机器标识符`r2 字节的进程 ID`r3字节的计数器，以随机值开始"

$bp1 = $bodyId.Paragraphs(1)
$bp1.Text = "默认创建"
$r = $bp1.InsertBefore("属于单字段索引的一种，")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("MongoDB")
$r.LanguageID = "en-US"
$bp1.LanguageID = "zh-CN"

$bp2 = $bodyId.Paragraphs(2)
$bp2.Text = "类型之一。"
$r = $bp2.InsertBefore("默认情况下，")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("_")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("id ")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("字段的类型为 ")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("ObjectID")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("，")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("是 ")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("MongoDB ")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("的 ")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("BSON ")
$r.LanguageID = "en-US"
$bp2.LanguageID = "zh-CN"

$bp3 = $bodyId.Paragraphs(3)
$bp3.Text = " 长度为 "
$r = $bp3.InsertBefore("ObjectID")
$r.LanguageID = "en-US"
$r = $r.InsertAfter($bp3.Text)
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("12 ")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("字节，由以下")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("4")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("个部分组成：")
$r.LanguageID = "zh-CN"
$bp3.Text = ""
$bp3.InsertBefore("ObjectID").LanguageID = "en-US"

$bp4 = $bodyId.Paragraphs(4)
$bp4.Text = "字节的"
$r = $bp4.InsertBefore("4 ")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("字节的")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("Unix")
$r.LanguageID = "en-US"
$r = $r.InsertAfter("时间戳")
$r.LanguageID = "zh-CN"
$r = $r.InsertAfter("，单位为秒")
$r.LanguageID = "zh-CN"
$bp4.ParagraphFormat.Bullet.Visible = 1
$bp4.ParagraphFormat.Bullet.Character = 108
$bp4.ParagraphFormat.Bullet.Font.Name = "Wingdings"

$bp5 = $bodyId.Paragraphs(5)
$bp5.Text = "字节的机器标识符"
$r = $bp5.InsertBefore("3 ")
$r.LanguageID = "en-US"
$bp5.LanguageID = "zh-CN"
$bp5.ParagraphFormat.Bullet.Visible = 1
$bp5.ParagraphFormat.Bullet.Character = 108
$bp5.ParagraphFormat.Bullet.Font.Name = "Wingdings"

$bp6 = $bodyId.Paragraphs(6)
$bp6.Text = "字节的进程 "
$r = $bp6.InsertBefore("2 ")
$r.LanguageID = "en-US"
$r = $bp6.InsertAfter("ID")
$r.LanguageID = "en-US"
$bp6.LanguageID = "zh-CN"
$bp6.ParagraphFormat.Bullet.Visible = 1
$bp6.ParagraphFormat.Bullet.Character = 108
$bp6.ParagraphFormat.Bullet.Font.Name = "Wingdings"

$bp7 = $bodyId.Paragraphs(7)
$bp7.Text = "字节的计数器，以随机值开始"
$r = $bp7.InsertBefore("3")
$r.LanguageID = "en-US"
$bp7.LanguageID = "zh-CN"
$bp7.ParagraphFormat.Bullet.Visible = 1
$bp7.ParagraphFormat.Bullet.Character = 108
$bp7.ParagraphFormat.Bullet.Font.Name = "Wingdings"

$nsId.Shapes.Item(2).TextFrame.TextRange.Paragraphs(3).ParagraphFormat.Bullet.Visible = 0

# --- New slide 5: "单字段索引" (title only, empty content) ---
$nsField = $p.Slides.Add(5, 2)
$nsField.Shapes.Item(1).Name = "标题 1"
$nsField.Shapes.Item(2).Name = "内容占位符 2"
$nsField.Shapes.Item(1).TextFrame.TextRange.Text = "单字段索引"
$nsField.Shapes.Item(1).TextFrame.TextRange.LanguageID = "zh-CN"

# --- New slide 6: "复合索引" (title only, empty content) ---
$nsComp = $p.Slides.Add(6, 2)
$nsComp.Shapes.Item(1).Name = "标题 1"
$nsComp.Shapes.Item(2).Name = "内容占位符 2"
$nsComp.Shapes.Item(1).TextFrame.TextRange.Text = "复合索引"
$nsComp.Shapes.Item(1).TextFrame.TextRange.LanguageID = "zh-CN"

# --- New slide 7: "多键索引" (title only, empty content) ---
$nsMulti = $p.Slides.Add(7, 2)
$nsMulti.Shapes.Item(1).Name = "标题 1"
$nsMulti.Shapes.Item(2).Name = "内容占位符 2"
$nsMulti.Shapes.Item(1).TextFrame.TextRange.Text = "多键索引"
$nsMulti.Shapes.Item(1).TextFrame.TextRange.LanguageID = "zh-CN"
